$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns W..BI (in order) on the header/comment row (row 15). The column
# "culture_collection" (W) is being removed; every comment to its right
# needs to shift one column to the left to mirror what Excel would do if
# the whole column (including its cell note) were deleted and everything
# past it slid over.
$cols = @("W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD","BE","BF","BG","BH","BI")

for ($i = 0; $i -lt $cols.Length - 1; $i++) {
    $dst = $ws.Range($cols[$i] + "15")
    $src = $ws.Range($cols[$i + 1] + "15")
    $text = $src.Comment.Text()
    $dst.Comment.Text($text) | Out-Null
}

# The final column's comment now has a duplicate one column to its left;
# drop it before the column cells themselves shift over.
$ws.Range($cols[$cols.Length - 1] + "15").Comment.Delete()

# Now remove the actual "culture_collection" column, sliding the cell
# values/shared strings left to match.
$ws.Columns("W").Delete()
